$d = $word.ActiveDocument

$replacements = @(
    @("734÷5=146, 4", "575÷9=63, 8"),
    @("352÷3=117, 1", "619÷2=309, 1"),
    @("275÷4=68, 3", "746÷2=373, 0"),
    @("946÷8=118, 2", "241÷5=48, 1"),
    @("216÷4=54, 0", "615÷7=87, 6"),
    @("917÷6=152, 5", "616÷3=205, 1"),
    @("912÷5=182, 2", "398÷8=49, 6"),
    @("719÷4=179, 3", "910÷3=303, 1"),
    @("659÷5=131, 4", "216÷9=24, 0"),
    @("366÷2=183, 0", "168÷6=28, 0"),
    @("350÷2=175, 0", "196÷6=32, 4"),
    @("857÷8=107, 1", "409÷8=51, 1"),
    @("247÷3=82, 1", "569÷3=189, 2"),
    @("909÷4=227, 1", "461÷3=153, 2"),
    @("443÷9=49, 2", "659÷9=73, 2"),
    @("250÷4=62, 2", "235÷2=117, 1"),
    @("437÷7=62, 3", "313÷6=52, 1"),
    @("106÷4=26, 2", "408÷2=204, 0"),
    @("376÷3=125, 1", "847÷5=169, 2"),
    @("797÷9=88, 5", "591÷5=118, 1"),
    @("430÷8=53, 6", "710÷7=101, 3"),
    @("266÷9=29, 5", "927÷7=132, 3"),
    @("650÷4=162, 2", "865÷9=96, 1"),
    @("495÷6=82, 3", "184÷9=20, 4"),
    @("678÷2=339, 0", "122÷8=15, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
